# stock update [kunai, jjk, dbz]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- JJK restock (rows 206-214, 216) ---
$ws.Range("B206").Value = 1
$ws.Range("B207").Value = 1
$ws.Range("B208").Value = 2
$ws.Range("B209").Value = 1
$ws.Range("B210").Value = 1
$ws.Range("B211").Value = 2
$ws.Range("B212").Value = 1
$ws.Range("B213").Value = 1
$ws.Range("B214").Value = 2
$ws.Range("B216").Value = 2

# --- DBZ restock (rows 256-268) ---
$ws.Range("B256").Value = 3
$ws.Range("B257").Value = 2
$ws.Range("B258").Value = 3
$ws.Range("B259").Value = 2
$ws.Range("B260").Value = 1
$ws.Range("B261").Value = 2
$ws.Range("B262").Value = 1
$ws.Range("B263").Value = 1
$ws.Range("B264").Value = 1
$ws.Range("B265").Value = 1
$ws.Range("B266").Value = 3
$ws.Range("B267").Value = 3
$ws.Range("B268").Value = 1

# --- New kunai products (rows 291-292) ---
# Shared-string insertion order matters: Hiraishin Kunai first, then Kunai,
# then kunai.jpg, then minato kunai.jpg (matches author's original edit order).
$ws.Range("A292").Value = "Hiraishin Kunai [Metal]"
$ws.Range("A291").Value = "Kunai [Metal]"
$ws.Range("D291").Value = "kunai.jpg"
$ws.Range("D292").Value = "minato kunai.jpg"

$ws.Range("E291").Value = "Others"
$ws.Range("E292").Value = "Others"

# Copy row formatting down from the last existing data row so the new rows
# look consistent with the rest of the table.
$ws.Range("A290:E290").Copy() | Out-Null
$ws.Range("A291:E291").PasteSpecial(-4122) | Out-Null
$ws.Range("A290:E290").Copy() | Out-Null
$ws.Range("A292:E292").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("B291").Value = 2
$ws.Range("C291").Value = 1200
$ws.Range("B292").Value = 2
$ws.Range("C292").Value = 1200

$ws.Range("D292").Select() | Out-Null
